# Recibo Nacional - Recibo combinado
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: updated order numbers ---
$ws.Range("A3").Value = 4500144688
$ws.Range("B3").Value = 7450043069886

# --- New column H width (~14.7109375 chars; engine quantizes to 1/6-char
#     pixel steps, so 13.8 is the closest input that lands on the nearest
#     representable width) ---
$ws.Columns("H").ColumnWidth = 13.8

# --- Row 20: new header cell H20 ---
$ws.Range("H20").Value = "Pedido Nacional Comb - METROMALL"

# --- Row 21: new cell H21 ---
$ws.Range("H21").Value = 4500144611

# --- Row 22: updated order + new combo-receipt header row ---
$ws.Range("A22").Value = 4500144691
$ws.Range("B22").Value = 2040000642174
$ws.Range("C22").Value = 380
$ws.Range("H22").Value = "EAN"
$ws.Range("I22").Value = "CANTIDAD"
$ws.Range("J22").Value = "EMPAQUE"
$ws.Range("K23").Value = "M005"
$ws.Range("K22").Value = "Centro PPK"
$ws.Range("L22").Value = "Cantidad"

# --- Row 23: new combo-receipt data row ---
$ws.Range("H23").Value = 7611501611195
$ws.Range("H23").NumberFormat = "0"
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = 12
$ws.Range("L23").Value = 12

# --- Row 24: new cell H24 ---
$ws.Range("H24").Value = 7822102021096
$ws.Range("H24").NumberFormat = "0"

# --- Row 25: new combo-receipt data row ---
$ws.Range("H25").Value = 7822102021102
$ws.Range("H25").NumberFormat = "0"

# --- Row 26: new order line + combo-receipt data row (highlighted) ---
$ws.Range("A26").Value = 4500144690
$ws.Range("B26").Value = 19048183262
$ws.Range("C26").Value = 50
$ws.Range("H26").Value = 7201902266188
$ws.Range("H26").NumberFormat = "0"
$ws.Range("H26").Interior.Color = 5296274
$ws.Range("I26").Value = 12

# --- Row 27: new combo-receipt data row (highlighted) ---
$ws.Range("H27").Value = 7201902266645
$ws.Range("H27").NumberFormat = "0"
$ws.Range("H27").Interior.Color = 5296274

# --- Row 28: new combo-receipt data row (highlighted) ---
$ws.Range("H28").Value = 7611501611126
$ws.Range("H28").NumberFormat = "0"
$ws.Range("H28").Interior.Color = 5296274

# --- View: scrolled + selection on H28 ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("H28").Select()
